$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 200007680
$ws.Range("I62").Value = 250007100
$ws.Range("K62").Value = 250007100
$ws.Range("M62").Value = -250006476
$ws.Range("H65").Value = 200007680
$ws.Range("I65").Value = 250007100
$ws.Range("K65").Value = 1250035500
$ws.Range("M65").Value = -1250032380
$ws.Range("H76").Value = 5348.385
$ws.Range("J76").Value = 5593.364
$ws.Range("L76").Value = 5593.364
$ws.Range("N76").Value = -6223.364
$ws.Range("H79").Value = 5348.385
$ws.Range("J79").Value = 5593.364
$ws.Range("L79").Value = 5593.364
$ws.Range("N79").Value = -7777.364
$ws.Range("H100").Value = 28000222
$ws.Range("I100").Value = 31440376
$ws.Range("J100").Value = 479000
$ws.Range("K100").Value = 31440376
$ws.Range("L100").Value = 479000
$ws.Range("M100").Value = -31439835
$ws.Range("N100").Value = -480082
$ws.Range("H132").Value = 1641779
$ws.Range("I132").Value = 2335.6853
$ws.Range("J132").Value = 14288913
$ws.Range("K132").Value = 7007.0559
$ws.Range("L132").Value = 42866739
$ws.Range("M132").Value = -4477.0559
$ws.Range("N132").Value = -42871799
$ws.Range("H137").Value = 530550.4
$ws.Range("I137").Value = 670697.3
$ws.Range("K137").Value = 2012091.9
$ws.Range("M137").Value = -2009541.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22707
$ws.Range("I32").Value = 17922.79
$ws.Range("J32").Value = 53007
$ws.Range("K32").Value = 17922.79
$ws.Range("L32").Value = 53007
$ws.Range("M32").Value = -17635.79
$ws.Range("N32").Value = -53581
$ws.Range("H45").Value = 88531.75
$ws.Range("I45").Value = 115682
$ws.Range("K45").Value = 115682
$ws.Range("M45").Value = -115305
$ws.Range("H63").Value = 1980.8
$ws.Range("I63").Value = 1976.25
$ws.Range("K63").Value = 1976.25
$ws.Range("M63").Value = -1290.25
$ws.Range("H66").Value = 1980.8
$ws.Range("I66").Value = 1976.25
$ws.Range("K66").Value = 9881.25
$ws.Range("M66").Value = -6449.25
$ws.Range("H96").Value = 74489
$ws.Range("J96").Value = 74489
$ws.Range("L96").Value = 74489
$ws.Range("N96").Value = -79981

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 36663.332
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3616
$ws.Range("H52").Value = 26498
$ws.Range("J52").Value = 26498
$ws.Range("L52").Value = 26498
$ws.Range("N52").Value = -27024
$ws.Range("H86").Value = 5861.522
$ws.Range("I86").Value = 6523.0557
$ws.Range("J86").Value = 3480
$ws.Range("K86").Value = 6523.0557
$ws.Range("L86").Value = 3480
$ws.Range("M86").Value = -5400.0557
$ws.Range("N86").Value = -5726
$ws.Range("H89").Value = 5861.522
$ws.Range("I89").Value = 6523.0557
$ws.Range("J89").Value = 3480
$ws.Range("K89").Value = 32615.2785
$ws.Range("L89").Value = 17400
$ws.Range("M89").Value = -26999.2785
$ws.Range("N89").Value = -28632
$ws.Range("H107").Value = 3050
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 3050
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 3050
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6890
$ws.Range("H121").Value = 26498
$ws.Range("J121").Value = 26498
$ws.Range("L121").Value = 26498
$ws.Range("N121").Value = -29992
$ws.Range("H134").Value = 4538.9287
$ws.Range("I134").Value = 3681.6365
$ws.Range("K134").Value = 11044.9095
$ws.Range("M134").Value = -8509.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3240.111
$ws.Range("I58").Value = 2191
$ws.Range("K58").Value = 2191
$ws.Range("M58").Value = -1988
$ws.Range("H99").Value = 6949898
$ws.Range("I99").Value = 15628909
$ws.Range("J99").Value = 6688.7
$ws.Range("K99").Value = 15628909
$ws.Range("L99").Value = 6688.7
$ws.Range("M99").Value = -15627411
$ws.Range("N99").Value = -9684.700000000001
$ws.Range("H107").Value = 7742.8
$ws.Range("I107").Value = 10277.454
$ws.Range("J107").Value = 772.5
$ws.Range("K107").Value = 10277.454
$ws.Range("L107").Value = 772.5
$ws.Range("M107").Value = -8357.454
$ws.Range("N107").Value = -4612.5
$ws.Range("H126").Value = 6949898
$ws.Range("I126").Value = 15628909
$ws.Range("J126").Value = 6688.7
$ws.Range("K126").Value = 46886727
$ws.Range("L126").Value = 20066.1
$ws.Range("M126").Value = -46884257
$ws.Range("N126").Value = -25006.1
$ws.Range("H132").Value = 6626.9
$ws.Range("I132").Value = 7372.8823
$ws.Range("K132").Value = 22118.6469
$ws.Range("M132").Value = -19588.6469
$ws.Range("H136").Value = 3240.111
$ws.Range("I136").Value = 2191
$ws.Range("K136").Value = 6573
$ws.Range("M136").Value = -4023
$ws.Range("H141").Value = 568931.5
$ws.Range("J141").Value = 598239.7
$ws.Range("L141").Value = 598239.7
$ws.Range("N141").Value = -608599.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 41678040
$ws.Range("J68").Value = 41678040
$ws.Range("L68").Value = 125034120
$ws.Range("N68").Value = -125035742
$ws.Range("H71").Value = 41678040
$ws.Range("J71").Value = 41678040
$ws.Range("L71").Value = 375102360
$ws.Range("N71").Value = -375110472
$ws.Range("H104").Value = 7058
$ws.Range("J104").Value = 6322.5
$ws.Range("L104").Value = 18967.5
$ws.Range("N104").Value = -24209.5
$ws.Range("H106").Value = 6523.75
$ws.Range("I106").Value = 3026
$ws.Range("J106").Value = 6841.727
$ws.Range("K106").Value = 9078
$ws.Range("L106").Value = 20525.181
$ws.Range("M106").Value = -8132
$ws.Range("N106").Value = -22417.181
$ws.Range("H131").Value = 5558.391
$ws.Range("I131").Value = 7174.6665
$ws.Range("J131").Value = 3795.182
$ws.Range("K131").Value = 21523.9995
$ws.Range("L131").Value = 11385.546
$ws.Range("M131").Value = -16483.9995
$ws.Range("N131").Value = -21465.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 321508500
$ws.Range("J98").Value = 321508500
$ws.Range("L98").Value = 321508500
$ws.Range("N98").Value = -321514490
$ws.Range("H102").Value = 3904.4285
$ws.Range("I102").Value = 4151.619
$ws.Range("J102").Value = 2421.2856
$ws.Range("K102").Value = 4151.619
$ws.Range("L102").Value = 2421.2856
$ws.Range("M102").Value = -2529.619
$ws.Range("N102").Value = -5665.2856
$ws.Range("H107").Value = 675.5357
$ws.Range("I107").Value = 560.05
$ws.Range("K107").Value = 560.05
$ws.Range("M107").Value = 1359.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 28330.334
$ws.Range("I7").Value = 38495.668
$ws.Range("K7").Value = 38495.668
$ws.Range("M7").Value = -38383.668
$ws.Range("H16").Value = 4534.524
$ws.Range("I16").Value = 4519.2354
$ws.Range("K16").Value = 4519.2354
$ws.Range("M16").Value = -4349.2354
$ws.Range("H40").Value = 86163.836
$ws.Range("J40").Value = 34994.668
$ws.Range("L40").Value = 34994.668
$ws.Range("N40").Value = -35266.668
$ws.Range("H46").Value = 1480.4
$ws.Range("J46").Value = 1971.5
$ws.Range("L46").Value = 1971.5
$ws.Range("N46").Value = -2347.5
$ws.Range("H122").Value = 8023.1875
$ws.Range("I122").Value = 9498.286
$ws.Range("K122").Value = 28494.858
$ws.Range("M122").Value = -26044.858
$ws.Range("H126").Value = 28330.334
$ws.Range("I126").Value = 38495.668
$ws.Range("K126").Value = 115487.004
$ws.Range("M126").Value = -113017.004
$ws.Range("H136").Value = 9014.15
$ws.Range("I136").Value = 10298.625
$ws.Range("K136").Value = 30895.875
$ws.Range("M136").Value = -28345.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 81687.836
$ws.Range("I74").Value = 4999
$ws.Range("J74").Value = 97025.60000000001
$ws.Range("K74").Value = 4999
$ws.Range("L74").Value = 97025.60000000001
$ws.Range("M74").Value = -4063
$ws.Range("N74").Value = -98897.60000000001
$ws.Range("H77").Value = 81687.836
$ws.Range("I77").Value = 4999
$ws.Range("J77").Value = 97025.60000000001
$ws.Range("K77").Value = 14997
$ws.Range("L77").Value = 291076.8
$ws.Range("M77").Value = -10317
$ws.Range("N77").Value = -300436.8
$ws.Range("H96").Value = 20002662
$ws.Range("I96").Value = 20002662
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 20002662
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -20001289
$ws.Range("N96").ClearContents()
$ws.Range("H136").Value = 815250.0600000001
$ws.Range("I136").Value = 1544345.8
$ws.Range("J136").Value = 5143.778
$ws.Range("K136").Value = 4633037.4
$ws.Range("L136").Value = 15431.334
$ws.Range("M136").Value = -4630487.4
$ws.Range("N136").Value = -20531.334
